# Generate Report for Archive
# - Update the "Ready for handoff" status text to "In Translation" everywhere
#   it appears (Overview sheet's per-language status cells, and the Status
#   column on each per-language sheet).
# - Narrow the "zh-cn"/"de-de" status columns (Overview col E/F, and col C on
#   each per-language sheet) from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Text change: "Ready for handoff" -> "In Translation" ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Column width changes: 17.2159881591797 -> 13.4101845877511 (stored width) ---
# The stored OOXML column width is derived from the character-based
# ColumnWidth through a pixel-grid rounding, so we pick the ColumnWidth value
# that rounds to the closest achievable stored width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
